# Add a "number formats" demo section to the styles example workbook,
# mirroring the existing "one row per style feature" pattern already in
# the sheet: column A gets a label, column B gets a sample numeric value
# (1.2) formatted with a distinct number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# label -> Excel NumberFormat string to apply to the adjacent B cell.
$formats = [ordered]@{
    "number"     = "0.00"
    "currency"   = "$#,##0.00"
    "accounting" = '_($* #,##0.00_);_($* (#,##0.00);_($* "-"??_);_(@_)'
    "short date" = "m/d/yyyy"
    "long date"  = "[$-x-sysdate]dddd, mmmm dd, yyyy"
    "time"       = "[$-x-systime]h:mm:ss AM/PM"
    "percentage" = "0.00%"
    "fraction"   = "# ?/?"
    "scientific" = "0.00E+00"
    "text"       = "@"
}

$row = 36
foreach ($label in $formats.Keys) {
    $ws.Cells.Item($row, 1).Value = $label
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = 1.2
    $cell.NumberFormat = $formats[$label]
    $row = $row + 1
}
